$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 was deleted (old "March paycheck #2" entry) which shifts the rows below
# it up by one; the row that becomes the new row 8 is a freshly added entry.
$ws.Rows("3").Delete()

$ws.Range("A3").Value = "April Paycheck #1"
$ws.Range("B3").Value = 1250
$ws.Range("C3").Value = 45017

$ws.Range("A8").Value = "January Paycheck #2"
$ws.Range("B8").Value = 1200
$ws.Range("C8").Value = 44941
$ws.Range("D8").Value = "Paycheck"
$ws.Range("C8").NumberFormat = $ws.Range("C7").NumberFormat
